$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"4.18e-007"
$ws.Range("E2").Value = [double]"2.97e-005"
$ws.Range("D3").Value = 0.051
$ws.Range("D4").Value = 0.09
$ws.Range("D5").Value = 0.11
$ws.Range("D6").Value = 0.128
$ws.Range("D7").Value = 0.172
$ws.Range("D8").Value = 0.175
$ws.Range("D9").Value = 0.228
$ws.Range("D10").Value = 0.326
$ws.Range("D11").Value = 0.331
$ws.Range("D12").Value = 0.337
$ws.Range("D13").Value = 0.417
$ws.Range("D14").Value = 0.448
$ws.Range("D15").Value = 0.455
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 0.498
$ws.Range("D17").Value = 0.519
$ws.Range("D18").Value = 0.533
$ws.Range("A19").Formula = '=HYPERLINK("pathways/Glycerol_Phosphate_Shuttle.csv","Glycerol Phosphate Shuttle")'
$ws.Range("B19").Value = 3
$ws.Range("D19").Value = 0.538
$ws.Range("A20").Formula = '=HYPERLINK("pathways/Beta_Alanine_Metabolism.csv","Beta-Alanine Metabolism")'
$ws.Range("B20").Value = 10
$ws.Range("D20").Value = 0.538
$ws.Range("D21").Value = 0.543
$ws.Range("D22").Value = 0.5590000000000001
$ws.Range("D23").Value = 0.569
$ws.Range("A24").Formula = '=HYPERLINK("pathways/Pyrimidine_Metabolism.csv","Pyrimidine Metabolism")'
$ws.Range("B24").Value = 13
$ws.Range("D24").Value = 0.64
$ws.Range("A25").Formula = '=HYPERLINK("pathways/Warburg_Effect.csv","Warburg Effect")'
$ws.Range("B25").Value = 18
$ws.Range("D25").Value = 0.64
$ws.Range("D26").Value = 0.652
$ws.Range("D27").Value = 0.654
$ws.Range("D28").Value = 0.698
$ws.Range("D29").Value = 0.698
$ws.Range("D30").Value = 0.698
$ws.Range("D31").Value = 0.71
$ws.Range("D32").Value = 0.712
$ws.Range("D33").Value = 0.716
$ws.Range("D34").Value = 0.734
$ws.Range("D35").Value = 0.734
$ws.Range("D37").Value = 0.742
$ws.Range("D38").Value = 0.783
$ws.Range("D39").Value = 0.799
$ws.Range("D40").Value = 0.8080000000000001
$ws.Range("A42").Formula = '=HYPERLINK("pathways/Gluconeogenesis.csv","Gluconeogenesis")'
$ws.Range("B42").Value = 10
$ws.Range("D42").Value = 0.833
$ws.Range("A43").Formula = '=HYPERLINK("pathways/Steroid_Biosynthesis.csv","Steroid Biosynthesis")'
$ws.Range("B43").Value = 3
$ws.Range("D43").Value = 0.835
$ws.Range("D44").Value = 0.843
$ws.Range("D47").Value = 0.853
$ws.Range("D48").Value = 0.855
$ws.Range("D50").Value = 0.864
$ws.Range("D51").Value = 0.866
$ws.Range("A52").Formula = '=HYPERLINK("pathways/Mitochondrial_Electron_Transport_Chain.csv","Mitochondrial Electron Transport Chain")'
$ws.Range("B52").Value = 8
$ws.Range("D52").Value = 0.887
$ws.Range("A53").Formula = '=HYPERLINK("pathways/Porphyrin_Metabolism.csv","Porphyrin Metabolism")'
$ws.Range("B53").Value = 4
$ws.Range("D53").Value = 0.888
$ws.Range("D54").Value = 0.899
$ws.Range("A55").Formula = '=HYPERLINK("pathways/Glutamate_Metabolism.csv","Glutamate Metabolism")'
$ws.Range("B55").Value = 14
$ws.Range("D55").Value = 0.9
$ws.Range("A56").Formula = '=HYPERLINK("pathways/Thyroid_hormone_synthesis.csv","Thyroid hormone synthesis")'
$ws.Range("B56").Value = 5
$ws.Range("D56").Value = 0.901
$ws.Range("A57").Formula = '=HYPERLINK("pathways/Alanine_Metabolism.csv","Alanine Metabolism")'
$ws.Range("B57").Value = 9
$ws.Range("D57").Value = 0.901
$ws.Range("A58").Formula = '=HYPERLINK("pathways/Lysine_Degradation.csv","Lysine Degradation")'
$ws.Range("B58").Value = 5
$ws.Range("D59").Value = 0.904
$ws.Range("D60").Value = 0.909
$ws.Range("D61").Value = 0.914
$ws.Range("D62").Value = 0.914
$ws.Range("D63").Value = 0.914
$ws.Range("D64").Value = 0.919
$ws.Range("D66").Value = 0.941

Write-Host "Applied all dins30_cameraPR updates"
